{"js": "// Update the date line and all 25 multiplication-practice answers in the\n// single table to the new values from the commit.\n//\n// Strategy: search the whole document body for each OLD string (exact,\n// case-sensitive, whole match) and replace it in place with the NEW\n// string. This is robust regardless of how the text happens to be\n// distributed across paragraphs / table cells, and does not depend on\n// assumptions about run-splitting.\n\nconst replacements = [\n  [\"2025-03-06 Thursday\", \"2025-03-07 Friday\"],\n  [\"919\u00d78=7352\", \"882\u00d75=4410\"],\n  [\"506\u00d79=4554\", \"560\u00d79=5040\"],\n  [\"528\u00d78=4224\", \"593\u00d79=5337\"],\n  [\"737\u00d79=6633\", \"823\u00d76=4938\"],\n  [\"254\u00d72=508\", \"634\u00d72=1268\"],\n  [\"453\u00d78=3624\", \"994\u00d79=8946\"],\n  [\"423\u00d72=846\", \"962\u00d78=7696\"],\n  [\"133\u00d73=399\", \"909\u00d74=3636\"],\n  [\"207\u00d73=621\", \"250\u00d76=1500\"],\n  [\"123\u00d75=615\", \"812\u00d79=7308\"],\n  [\"157\u00d75=785\", \"829\u00d77=5803\"],\n  [\"323\u00d75=1615\", \"591\u00d74=2364\"],\n  [\"533\u00d73=1599\", \"745\u00d77=5215\"],\n  [\"826\u00d73=2478\", \"150\u00d78=1200\"],\n  [\"481\u00d78=3848\", \"766\u00d79=6894\"],\n  [\"303\u00d74=1212\", \"821\u00d74=3284\"],\n  [\"325\u00d76=1950\", \"931\u00d72=1862\"],\n  [\"508\u00d76=3048\", \"471\u00d76=2826\"],\n  [\"463\u00d76=2778\", \"691\u00d79=6219\"],\n  [\"165\u00d74=660\", \"182\u00d78=1456\"],\n  [\"207\u00d74=828\", \"924\u00d75=4620\"],\n  [\"968\u00d74=3872\", \"457\u00d73=1371\"],\n  [\"500\u00d77=3500\", \"407\u00d73=1221\"],\n  [\"357\u00d79=3213\", \"471\u00d78=3768\"],\n  [\"868\u00d75=4340\", \"426\u00d73=1278\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 multiplication-practice answers in the\n# single table to the new values from the commit, using Find/Replace across\n# the whole document story (Content range covers the body + the table).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-06 Thursday\", \"2025-03-07 Friday\"),\n    @(\"919\u00d78=7352\", \"882\u00d75=4410\"),\n    @(\"506\u00d79=4554\", \"560\u00d79=5040\"),\n    @(\"528\u00d78=4224\", \"593\u00d79=5337\"),\n    @(\"737\u00d79=6633\", \"823\u00d76=4938\"),\n    @(\"254\u00d72=508\", \"634\u00d72=1268\"),\n    @(\"453\u00d78=3624\", \"994\u00d79=8946\"),\n    @(\"423\u00d72=846\", \"962\u00d78=7696\"),\n    @(\"133\u00d73=399\", \"909\u00d74=3636\"),\n    @(\"207\u00d73=621\", \"250\u00d76=1500\"),\n    @(\"123\u00d75=615\", \"812\u00d79=7308\"),\n    @(\"157\u00d75=785\", \"829\u00d77=5803\"),\n    @(\"323\u00d75=1615\", \"591\u00d74=2364\"),\n    @(\"533\u00d73=1599\", \"745\u00d77=5215\"),\n    @(\"826\u00d73=2478\", \"150\u00d78=1200\"),\n    @(\"481\u00d78=3848\", \"766\u00d79=6894\"),\n    @(\"303\u00d74=1212\", \"821\u00d74=3284\"),\n    @(\"325\u00d76=1950\", \"931\u00d72=1862\"),\n    @(\"508\u00d76=3048\", \"471\u00d76=2826\"),\n    @(\"463\u00d76=2778\", \"691\u00d79=6219\"),\n    @(\"165\u00d74=660\", \"182\u00d78=1456\"),\n    @(\"207\u00d74=828\", \"924\u00d75=4620\"),\n    @(\"968\u00d74=3872\", \"457\u00d73=1371\"),\n    @(\"500\u00d77=3500\", \"407\u00d73=1221\"),\n    @(\"357\u00d79=3213\", \"471\u00d78=3768\"),\n    @(\"868\u00d75=4340\", \"426\u00d73=1278\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
